# smpte.xlsx edit: add support for SMPTE 4x3 row-14 based ratios (rows 16-19)
# and a new EBU test-pattern block (rows 21-26), plus new bar colors
# (I25:I31) and move the active-cell selection to I32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 16-19: re-point the %Height/%Width formulas at the SMPTE 4x3
#     header row (row 14, B14/C14) instead of the SMPTE 16x9 header (row 2).
$ws.Range("D16").Formula = '=B16/$B$14'
$ws.Range("E16").Formula = '=C16/$C$14'
$ws.Range("D17").Formula = '=B17/$B$14'
$ws.Range("E17").Formula = '=C17/$C$14'
$ws.Range("D18").Formula = '=B18/$B$14'
$ws.Range("E18").Formula = '=C18/$C$14'
$ws.Range("D19").Formula = '=B19/$B$14'
$ws.Range("E19").Formula = '=C19/$C$14'

# --- Row 20: blank separator row - clear any formula/value, keep formatting.
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()

# --- Row 21: new EBU header row (like rows 2 and 14) holding the base
#     height/width for the new EBU test pattern; no ratio formulas here.
$ws.Range("B21").Value = 150
$ws.Range("C21").Value = 200
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()

# --- Row 22: first EBU data row, with its own (non-shared) ratio formulas.
$ws.Range("B22").Value = 150
$ws.Range("C22").Value = 25
$ws.Range("D22").Formula = '=B22/$B$21'
$ws.Range("E22").Formula = '=C22/$C$21'

# --- Rows 23-26: remaining EBU data rows, ratios against the new row 21
#     header (previously these referenced row 2's header).
$ws.Range("D23").Formula = '=B23/$B$21'
$ws.Range("E23").Formula = '=C23/$C$21'
$ws.Range("D24").Formula = '=B24/$B$21'
$ws.Range("E24").Formula = '=C24/$C$21'
$ws.Range("D25").Formula = '=B25/$B$21'
$ws.Range("E25").Formula = '=C25/$C$21'
$ws.Range("D26").Formula = '=B26/$B$21'
$ws.Range("E26").Formula = '=C26/$C$21'

# --- Column I: new EBU bar colors (previously blank cells), sourced from
#     the newly added shared strings.
$ws.Range("I25").Value = "BFBFBF"
$ws.Range("I26").Value = "BFBF00"
$ws.Range("I27").Value = "00BFBF"
$ws.Range("I28").Value = "00BF00"
$ws.Range("I29").Value = "BF0000"
$ws.Range("I30").Value = "0000BF"

# --- Row 31: new trailing marker cell.
$ws.Range("I31").Value = 0

# --- Move the active cell / selection down to I32 to match the saved view.
[void]$ws.Range("I32").Select()
